$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM rows (10-20)
$ws.Range("E10").Value = "AndersonPP"
$ws.Range("H10").Value = 3

$ws.Range("E11").Value = "Zener Diode"
$ws.Range("H11").Value = 2

$ws.Range("E12").Value = "510 Resistor"
$ws.Range("F12").Value = "ESR03EZPJ511"
$ws.Range("G12").Value = "Digi-key"
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = 0.1
$ws.Range("J12").Value = "https://www.digikey.com/en/products/detail/rohm-semiconductor/ESR03EZPJ511/1762937"

$ws.Range("E13").Value = "330 Resistor"
$ws.Range("F13").Value = "ESR03EZPJ331"
$ws.Range("G13").Value = "Digi-key"
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 0.1
$ws.Range("J13").Value = "https://www.digikey.com/en/products/detail/rohm-semiconductor/ESR03EZPJ331/1762730"

$ws.Range("E14").Value = "100 Resistor "
$ws.Range("F14").Value = "ESR03EZPJ101"
$ws.Range("G14").Value = "Digi-key"
$ws.Range("H14").Value = 9
$ws.Range("I14").Value = 0.1
$ws.Range("J14").Value = "https://www.digikey.com/en/products/detail/rohm-semiconductor/ESR03EZPJ101/1983452"

$ws.Range("E15").Value = "43 Resistor "
$ws.Range("F15").Value = "ESR10EZPJ430"
$ws.Range("G15").Value = "Digi-key"
$ws.Range("H15").Value = 11
$ws.Range("I15").Value = 0.1
$ws.Range("J15").Value = "https://www.digikey.com/en/products/detail/rohm-semiconductor/ESR10EZPJ430/1762819"

$ws.Range("E16").Value = "0.1uF Capacitor"
$ws.Range("F16").Value = "CL10B104KO8NNNC"
$ws.Range("G16").Value = "Digi-key "
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = "https://www.digikey.com/en/products/detail/samsung-electro-mechanics/CL10B104KO8NNNC/3889091"

$ws.Range("E17").Value = "Teensy "
$ws.Range("H17").Value = 1

$ws.Range("E18").Value = "Molex Connectors"
$ws.Range("H18").Value = 8

$ws.Range("E19").Value = "OKI 5V"
$ws.Range("F19").Value = "OKI-78SR-5/1.5-W36H-C"
$ws.Range("G19").Value = "Digi-key"
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 4.3
$ws.Range("J19").Value = "https://www.digikey.com/en/products/detail/murata-power-solutions-inc/OKI-78SR-5-1.5-W36H-C/3438675"

$ws.Range("E20").Value = "OKI 3.3V"
$ws.Range("F20").Value = "OKI-78SR-3.3/1.5-W36H-C"
$ws.Range("G20").Value = "Digi-key"
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 4.3
$ws.Range("J20").Value = "https://www.digikey.com/en/products/detail/murata-power-solutions-inc/OKI-78SR-3.3-1.5-W36H-C/4878851"

# Restore the selection state that was present when the file was saved
$ws.Range("B8").Select() | Out-Null
